# Updated symbol list on Sun Jan  8 13:17:34 UTC 2023 with GitHub Actions
# Refreshes Price (D), Volume(1h) (E) and Hora (G) columns for the crypto table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "261.39"),
    @("E2", "0.36%"),
    @("G2", "13"),
    @("D3", "26.64"),
    @("E3", "-2.46%"),
    @("G3", "13"),
    @("D4", "4.702"),
    @("E4", "0.21%"),
    @("G4", "13"),
    @("D5", "0.06079"),
    @("E5", "-0.52%"),
    @("G5", "13"),
    @("D6", "6.705"),
    @("E6", "0.71%"),
    @("G6", "13"),
    @("D7", "0.8511"),
    @("E7", "-0.17%"),
    @("G7", "13"),
    @("D8", "0.9142"),
    @("E8", "-0.65%"),
    @("G8", "13"),
    @("D9", "0.1403"),
    @("E9", "0.17%"),
    @("G9", "13"),
    @("D10", "0.05185"),
    @("E10", "9.93%"),
    @("G10", "13"),
    @("E11", "0.01%"),
    @("G11", "13"),
    @("D12", "0.03117"),
    @("E12", "1.41%"),
    @("G12", "13"),
    @("D13", "0.09043"),
    @("E13", "-0.15%"),
    @("G13", "13"),
    @("D14", "0.001533"),
    @("E14", "-0.25%"),
    @("G14", "13"),
    @("D15", "0.0006165"),
    @("E15", "1.06%"),
    @("G15", "13"),
    @("D16", "0.006108"),
    @("E16", "1.51%"),
    @("G16", "13"),
    @("D17", "3.450"),
    @("E17", "-0.09%"),
    @("G17", "13"),
    @("E18", "0.82%"),
    @("G18", "13"),
    @("D19", "2.166"),
    @("E19", "0.14%"),
    @("G19", "13"),
    @("G20", "13"),
    @("D21", "0.1300"),
    @("E21", "-0.35%"),
    @("G21", "13"),
    @("D22", "4.097"),
    @("E22", "-0.07%"),
    @("G22", "13"),
    @("D23", "0.04234"),
    @("E23", "0.22%"),
    @("G23", "13"),
    @("E24", "-3.50%"),
    @("G24", "13"),
    @("D25", "0.004064"),
    @("E25", "6.77%"),
    @("G25", "13"),
    @("E26", "0.04%"),
    @("G26", "13"),
    @("E27", "4.14%"),
    @("G27", "13"),
    @("G28", "13"),
    @("G29", "13"),
    @("G30", "13"),
    @("G31", "13"),
    @("G32", "13"),
    @("G33", "13"),
    @("G34", "13"),
    @("G35", "13"),
    @("G36", "13"),
    @("G37", "13"),
    @("G38", "13"),
    @("G39", "13"),
    @("D40", "0.03950"),
    @("E40", "2.48%"),
    @("G40", "13"),
    @("E41", "-0.25%"),
    @("G41", "13"),
    @("D42", "0.004176"),
    @("E42", "1.94%"),
    @("G42", "13"),
    @("D43", "0.01389"),
    @("E43", "-14.85%"),
    @("G43", "13"),
    @("D44", "0.002066"),
    @("E44", "-6.79%"),
    @("G44", "13"),
    @("D45", "0.00005115"),
    @("E45", "-0.83%"),
    @("G45", "13"),
    @("E46", "0.06%"),
    @("G46", "13"),
    @("D47", "0.02121"),
    @("E47", "-60.72%"),
    @("G47", "13"),
    @("D48", "0.2584"),
    @("E48", "90.70%"),
    @("G48", "13"),
    @("E49", "0.06%"),
    @("G49", "13"),
    @("D50", "0.0002001"),
    @("E50", "0.06%"),
    @("G50", "13"),
    @("G51", "13")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    # Force text storage so values like "261.39" / "13" stay strings
    # (matching the inline-string cells already used throughout the sheet)
    # instead of being auto-converted to numbers by Excel.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}
